$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 28 (pushes existing rows 28:63 down to 29:64)
$ws.Rows("28:28").Insert()

# Fill the new row 28 with the new weekly price-record data.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R mirror the surrounding Achicoria / Quillota
# records; D,J,K,L,M,P carry the new observation's values.
$ws.Range("A28").Value = 9
$ws.Range("B28").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C28").Value = "Metropolitana"
$ws.Range("D28").Value = 45079
$ws.Range("E28").Value = 13
$ws.Range("F28").Value = 100112010
$ws.Range("G28").Value = "Achicoria"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 160
$ws.Range("K28").Value = 7000
$ws.Range("L28").Value = 7000
$ws.Range("M28").Value = 7000
$ws.Range("N28").Value = "$/caja 16 unidades"
$ws.Range("O28").Value = "Provincia de Quillota"
$ws.Range("P28").Value = 438
$ws.Range("Q28").Value = 16
$ws.Range("R28").Value = "Hortaliza"
